$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update execution times for the two remaining test cases
$ws.Range("E2").Value = "1845 ms"
$ws.Range("E3").Value = "2389 ms"

# Remove the third test case row entirely (TC_DeleteADM_03)
$ws.Rows(4).Delete()
